# Update the Odds_* columns (G:BD) for the rematched FlashScore rows (2, 3, 4, 7)
# on Sheet1, per the latest odds feed snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 1.83
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 2.6
$ws.Range("L2").Value = 5
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 7.5
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 19
$ws.Range("AC2").Value = 7
$ws.Range("AD2").Value = 6.5
$ws.Range("AE2").Value = 19
$ws.Range("AG2").Value = 10
$ws.Range("AH2").Value = 23
$ws.Range("AI2").Value = 17
$ws.Range("AJ2").Value = 51
$ws.Range("AK2").Value = 41
$ws.Range("AL2").Value = 51
$ws.Range("AN2").Value = 3.6
$ws.Range("AO2").Value = 10
$ws.Range("AS2").Value = 251
$ws.Range("AU2").Value = 9.5
$ws.Range("AW2").Value = 6
$ws.Range("AX2").Value = 29
$ws.Range("AY2").Value = 41
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 151
$ws.Range("BB2").Value = 351

# Row 3
$ws.Range("Q3").Value = 1.98
$ws.Range("R3").Value = 1.88

# Row 4
$ws.Range("H4").Value = 2.92
$ws.Range("I4").Value = 2.92
$ws.Range("K4").Value = 2.02
$ws.Range("L4").Value = 3.55
$ws.Range("M4").Value = 1.09
$ws.Range("N4").Value = 6.2
$ws.Range("O4").Value = 1.39
$ws.Range("P4").Value = 2.77
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.65
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.62
$ws.Range("U4").Value = 1.82
$ws.Range("V4").Value = 1.88
$ws.Range("W4").Value = 7.3
$ws.Range("Z4").Value = 27
$ws.Range("AA4").Value = 21
$ws.Range("AB4").Value = 32
$ws.Range("AC4").Value = 6.2
$ws.Range("AD4").Value = 5.8
$ws.Range("AE4").Value = 14
$ws.Range("AF4").Value = 75
$ws.Range("AG4").Value = 8
$ws.Range("AI4").Value = 10.75
$ws.Range("AK4").Value = 28
$ws.Range("AL4").Value = 37
$ws.Range("AM4").Value = 600
$ws.Range("AR4").Value = 90
$ws.Range("AT4").Value = 2.62
$ws.Range("AU4").Value = 6.8
$ws.Range("AX4").Value = 16.5
$ws.Range("AZ4").Value = 80

# Row 7
$ws.Range("G7").Value = 4.7
$ws.Range("H7").Value = 3.5
$ws.Range("I7").Value = 1.65
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 2.12
$ws.Range("L7").Value = 2.27
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 7.1
$ws.Range("O7").Value = 1.32
$ws.Range("P7").Value = 3.1
$ws.Range("Q7").Value = 1.95
$ws.Range("R7").Value = 1.78
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.67
$ws.Range("U7").Value = 1.91
$ws.Range("V7").Value = 1.8
$ws.Range("W7").Value = 12.5
$ws.Range("X7").Value = 27
$ws.Range("Y7").Value = 15.5
$ws.Range("AB7").Value = 55
$ws.Range("AC7").Value = 7.1
$ws.Range("AD7").Value = 7
$ws.Range("AE7").Value = 17
$ws.Range("AF7").Value = 90
$ws.Range("AG7").Value = 6.3
$ws.Range("AH7").Value = 7.4
$ws.Range("AI7").Value = 8.25
$ws.Range("AJ7").Value = 12.5
$ws.Range("AK7").Value = 13.5
$ws.Range("AL7").Value = 29
$ws.Range("AM7").Value = 700
$ws.Range("AN7").Value = 6.4
$ws.Range("AO7").Value = 28
$ws.Range("AS7").Value = 500
$ws.Range("AT7").Value = 2.67
$ws.Range("AU7").Value = 7.8
$ws.Range("AV7").Value = 80
$ws.Range("AX7").Value = 8.25
$ws.Range("AY7").Value = 19
$ws.Range("AZ7").Value = 29
$ws.Range("BA7").Value = 65
$ws.Range("BB7").Value = 300
